$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 8
$ws.Range("F7").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("F17").Value = -1
$ws.Range("F19").Value = -7
$ws.Range("F24").Value = -6
